$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price (column D) value is numeric-looking and must be
# forced to remain text (matching the source data, which stores prices as
# inline strings, not numbers) by pre-formatting the cell as Text before
# assigning the value - otherwise Excel auto-converts e.g. "26.30" -> 26.3.

$ws.Range("D2").Value = '58.144.36'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '3.127.70'
$ws.Range("E3").Value = '  +1.48%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("E5").Value = '  +1.53%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '142.59'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.124.94'
$ws.Range("E8").Value = '  +1.32%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.446'
$ws.Range("E9").Value = '  +1.57%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '7.15'
$ws.Range("E10").Value = '  -2.31%  '
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("E12").Value = '  +2.30%  '
$ws.Range("D13").Value = '3.662.50'
$ws.Range("E13").Value = '  +1.40%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.135'
$ws.Range("E14").Value = '  +3.55%  '
$ws.Range("E15").Value = '  -3.56%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.0000165'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '58.188.05'
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("D18").Value = '3.130.65'
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("E19").Value = '  -0.26%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '12.81'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("E21").Value = '  -1.52%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '342.55'
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("E23").Value = '  +0.02%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.513'
$ws.Range("E24").Value = '  +1.94%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '67.64'
$ws.Range("E25").Value = '  +2.97%  '
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("D28").Value = '0.0₃0930'
$ws.Range("E28").Value = '  +1.84%  '
$ws.Range("E29").Value = '  -0.01%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '7.36'
$ws.Range("E30").Value = '  +2.18%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '6.39'
$ws.Range("E31").Value = '  -2.56%  '
$ws.Range("E32").Value = '  +2.20%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '21.11'
$ws.Range("E33").Value = '  +0.70%  '
$ws.Range("E34").Value = '  -1.15%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '158.47'
$ws.Range("E35").Value = '  +2.83%  '
$ws.Range("E36").Value = '  +3.74%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '6.19'
$ws.Range("E37").Value = '  +2.46%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '26.30'
$ws.Range("E38").Value = '  -2.24%  '
$ws.Range("E39").Value = '  -4.06%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.66'
$ws.Range("E40").Value = '  +11.91%  '
$ws.Range("E41").Value = '  -1.53%  '
$ws.Range("E42").Value = '  +4.79%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.99'
$ws.Range("E43").Value = '  +2.68%  '
$ws.Range("D44").Value = '3.166.09'
$ws.Range("E44").Value = '  +1.28%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '36.65'
$ws.Range("E45").Value = '  -0.16%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.999'
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("E47").Value = '  +2.77%  '
$ws.Range("D48").Value = '2.257.08'
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("E49").Value = '  +4.90%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '6.10'
$ws.Range("E50").Value = '  +1.99%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '20.64'
$ws.Range("E51").Value = '  -0.05%  '
